# Regenerate the lattice-multiplication practice problems in the single
# 5x3 table so that the document matches the new exercise set while
# preserving the existing run formatting (sz=32 run props, <w:br/> line
# breaks inside a single run per cell).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

# Each entry is: row, col, top line ("A x B"), second line
# ("  d1    d2"), third row-label line, fourth row-label line.
# These are joined with vertical-tab (0x0B) characters, which Word's
# COM object model uses to represent manual line breaks (<w:br/>) when
# read from / written to Range.Text.
$cells = @(
    @{ Row=1; Col=1; Lines=@("17 x 67", "  6    7", "  ----", "1|    |", "7|    |") },
    @{ Row=1; Col=2; Lines=@("94 x 23", "  2    3", "  ----", "9|    |", "4|    |") },
    @{ Row=1; Col=3; Lines=@("94 x 57", "  5    7", "  ----", "9|    |", "4|    |") },

    @{ Row=2; Col=1; Lines=@("19 x 18", "  1    8", "  ----", "1|    |", "9|    |") },
    @{ Row=2; Col=2; Lines=@("53 x 78", "  7    8", "  ----", "5|    |", "3|    |") },
    @{ Row=2; Col=3; Lines=@("90 x 90", "  9    0", "  ----", "9|    |", "0|    |") },

    @{ Row=3; Col=1; Lines=@("42 x 34", "  3    4", "  ----", "4|    |", "2|    |") },
    @{ Row=3; Col=2; Lines=@("25 x 60", "  6    0", "  ----", "2|    |", "5|    |") },
    @{ Row=3; Col=3; Lines=@("17 x 53", "  5    3", "  ----", "1|    |", "7|    |") },

    @{ Row=4; Col=1; Lines=@("69 x 50", "  5    0", "  ----", "6|    |", "9|    |") },
    @{ Row=4; Col=2; Lines=@("91 x 84", "  8    4", "  ----", "9|    |", "1|    |") },
    @{ Row=4; Col=3; Lines=@("78 x 32", "  3    2", "  ----", "7|    |", "8|    |") },

    @{ Row=5; Col=1; Lines=@("25 x 46", "  4    6", "  ----", "2|    |", "5|    |") },
    @{ Row=5; Col=2; Lines=@("79 x 38", "  3    8", "  ----", "7|    |", "9|    |") },
    @{ Row=5; Col=3; Lines=@("19 x 93", "  9    3", "  ----", "1|    |", "9|    |") }
)

foreach ($entry in $cells) {
    $cell = $t.Cell($entry.Row, $entry.Col)
    $newText = [string]::Join($vtab, $entry.Lines)
    $cell.Range.Text = $newText
}
